# Gantry Parts List update
# - Adds a second "V-Slot Linear Rail 20mm x 60mm" line (1500mm length) in row 11
# - Adds a new "90 Degree Joining Plate" line in row 12
# - Adds a note "u.;" in E15
# - Extends the Total column formula/SUM to cover the new rows
# - Updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 11: V-Slot Linear Rail 20mm x 60mm, 1500mm, Silver
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "V-Slot Linear Rail 20mm x 60mm"
$ws.Range("C11").Value = "http://openbuildspartstore.com/v-slot-linear-rail/"

# ---------------------------------------------------------------------------
# Row 12: 90 Degree Joining Plate (entered before B11 so the shared-string
# table is populated in the same order as the source workbook)
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "http://openbuildspartstore.com/90-degree-joining-plate/"
$ws.Range("A12").Value = "90 Degree Joining Plate"
$ws.Range("B12").Value = "Connect X axis to Y axis"

$ws.Range("B11").Value = "1500mm, Silver"

# ---------------------------------------------------------------------------
# Row 15: stray note
# ---------------------------------------------------------------------------
$ws.Range("E15").Value = "u.;"

# ---------------------------------------------------------------------------
# Hyperlinks for the two new URL cells (reuse the look of the existing
# hyperlink cells by pasting the format from C10 after the link is added)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C11"), "http://openbuildspartstore.com/v-slot-linear-rail/")
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C12"), "http://openbuildspartstore.com/90-degree-joining-plate/")
$ws.Range("C10").Copy()
$ws.Range("C12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Quantities / prices / line totals
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 23
$ws.Range("F11").Formula = "=D11*E11"

$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 4.8
$ws.Range("F12").Formula = "=D12*E12"

# ---------------------------------------------------------------------------
# Grand total now covers the new rows
# ---------------------------------------------------------------------------
$ws.Range("F19").Formula = "=SUM(F2:F13)"

# ---------------------------------------------------------------------------
# View state: selection moves to E15, with the window scrolled so row 4 /
# column B is the top-left visible cell
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E15").Select()
